# 2021 CBOC sign-in sheet — "16-End" tab: add the MON 5/31 (Memorial Day,
# federal holiday) column as a highlighted/closed column, matching the
# other weekend/holiday column pairs (it mirrors the adjacent SUN 5/30
# column AD/AE in both formatting and the "X"/"X" closed markers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "16-End"

# --- Column widths: AF (32) / AG (33) go from 4.5 -> 2.5, matching the
# other narrow weekend columns (copy the width straight from column AE/31). ---
$refWidth = $ws.Columns.Item(31).ColumnWidth
$ws.Columns.Item(32).ColumnWidth = $refWidth
$ws.Columns.Item(33).ColumnWidth = $refWidth

# --- Row-by-row: copy the formatting of the AD/AE (SUN) pair onto the
# AF/AG (MON 5/31) pair for every data row, then set the matching content. ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("AD$r").Copy()
    $ws.Range("AF$r").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("AE$r").Copy()
    $ws.Range("AG$r").PasteSpecial(-4122)   # xlPasteFormats
}

# Row 2: day-of-week header ("MON" stays the same text, just re-styled)
$ws.Range("AF2").Value = "MON"
$ws.Range("AG2").Value = ""

# Row 3: date header (31 stays the same value, just re-styled)
$ws.Range("AF3").Value = 31
$ws.Range("AG3").Value = ""

# Row 4: "Tech" / "Time of Arrival" column captions (unchanged text, re-styled)
$ws.Range("AF4").Value = "Tech"
$ws.Range("AG4").Value = "Time of Arrival"

# Rows 5,8,11,14,17,20,23,26 -> CBOC/CORE location name rows: mark closed ("X")
foreach ($r in 5,8,11,14,17,20,23,26) {
    $ws.Range("AF$r").Value = "X"
    $ws.Range("AG$r").Value = "X"
}

# Rows 6,9,12,15,18,21,24 -> "Frozen" rows: mark closed ("X")
foreach ($r in 6,9,12,15,18,21,24) {
    $ws.Range("AF$r").Value = "X"
    $ws.Range("AG$r").Value = "X"
}

# Rows 7,10,13,16,19,22,25 -> blank spacer rows: stay empty (formatting only)
foreach ($r in 7,10,13,16,19,22,25) {
    $ws.Range("AF$r").Value = ""
    $ws.Range("AG$r").Value = ""
}

# Row 27 -> last "Frozen" row: mark closed ("X")
$ws.Range("AF27").Value = "X"
$ws.Range("AG27").Value = "X"
